# Apply scheduled market-data refresh to the per-job Leve profit tables.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) has a Table_<job>
# covering A1:N141 with current/average market prices (H:N) refreshed here.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9249.333000000001
$ws.Range("J32").Value = 9374.5
$ws.Range("L32").Value = 9374.5
$ws.Range("N32").Value = -10026.5
$ws.Range("H40").Value = 62582028
$ws.Range("I40").Value = 99499.5
$ws.Range("J40").Value = 83409540
$ws.Range("K40").Value = 99499.5
$ws.Range("L40").Value = 83409540
$ws.Range("M40").Value = -99324.5
$ws.Range("N40").Value = -83409890
$ws.Range("H41").Value = 3065.4
$ws.Range("J41").Value = 3124.125
$ws.Range("L41").Value = 3124.125
$ws.Range("N41").Value = -4004.125
$ws.Range("H64").Value = 14784030
$ws.Range("J64").Value = 29415868
$ws.Range("L64").Value = 29415868
$ws.Range("N64").Value = -29416364
$ws.Range("H67").Value = 14784030
$ws.Range("J67").Value = 29415868
$ws.Range("L67").Value = 29415868
$ws.Range("N67").Value = -29417584
$ws.Range("H94").Value = 13897499
$ws.Range("I94").Value = 18523334
$ws.Range("K94").Value = 18523334
$ws.Range("M94").Value = -18522883
$ws.Range("H100").Value = 1443.5
$ws.Range("I100").Value = 1443.5
$ws.Range("K100").Value = 1443.5
$ws.Range("M100").Value = -902.5
$ws.Range("H132").Value = 11564.702
$ws.Range("I132").Value = 4675.5625
$ws.Range("K132").Value = 14026.6875
$ws.Range("M132").Value = -11496.6875
$ws.Range("H141").Value = 6498.476
$ws.Range("I141").Value = 7075.8887
$ws.Range("K141").Value = 21227.6661
$ws.Range("M141").Value = -16047.6661

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4922.86
$ws.Range("I32").Value = 2492.7778
$ws.Range("J32").Value = 14035.667
$ws.Range("K32").Value = 2492.7778
$ws.Range("L32").Value = 14035.667
$ws.Range("M32").Value = -2205.7778
$ws.Range("N32").Value = -14609.667
$ws.Range("H45").Value = 3158.842
$ws.Range("I45").Value = 2901.5334
$ws.Range("K45").Value = 2901.5334
$ws.Range("M45").Value = -2524.5334
$ws.Range("H97").Value = 813.94446
$ws.Range("I97").Value = 975.0714
$ws.Range("K97").Value = 975.0714
$ws.Range("M97").Value = -479.0714
$ws.Range("H102").Value = 362215.47
$ws.Range("I102").Value = 686091.9
$ws.Range("K102").Value = 686091.9
$ws.Range("M102").Value = -684469.9
$ws.Range("H122").Value = 6264.5557
$ws.Range("I122").Value = 2445.5
$ws.Range("J122").Value = 9319.799999999999
$ws.Range("K122").Value = 7336.5
$ws.Range("L122").Value = 27959.4
$ws.Range("M122").Value = -4886.5
$ws.Range("N122").Value = -32859.39999999999
$ws.Range("H132").Value = 18936.615
$ws.Range("I132").Value = 23262.77
$ws.Range("K132").Value = 69788.31
$ws.Range("M132").Value = -67258.31

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2363.375
$ws.Range("I86").Value = 2266.5
$ws.Range("J86").Value = 2654
$ws.Range("K86").Value = 2266.5
$ws.Range("L86").Value = 2654
$ws.Range("M86").Value = -1143.5
$ws.Range("N86").Value = -4900
$ws.Range("H89").Value = 2363.375
$ws.Range("I89").Value = 2266.5
$ws.Range("J89").Value = 2654
$ws.Range("K89").Value = 11332.5
$ws.Range("L89").Value = 13270
$ws.Range("M89").Value = -5716.5
$ws.Range("N89").Value = -24502
$ws.Range("H99").Value = 1226800.1
$ws.Range("I99").Value = 1737140.1
$ws.Range("J99").Value = 1984.2
$ws.Range("K99").Value = 1737140.1
$ws.Range("L99").Value = 1984.2
$ws.Range("M99").Value = -1735642.1
$ws.Range("N99").Value = -4980.2
$ws.Range("H134").Value = 3555.9048
$ws.Range("I134").Value = 1559.7693
$ws.Range("J134").Value = 6799.625
$ws.Range("K134").Value = 4679.3079
$ws.Range("L134").Value = 20398.875
$ws.Range("M134").Value = -2144.3079
$ws.Range("N134").Value = -25468.875
$ws.Range("H139").Value = 99540
$ws.Range("J139").Value = 99540
$ws.Range("L139").Value = 99540
$ws.Range("N139").Value = -109820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7204
$ws.Range("I31").Value = 2399.5
$ws.Range("J31").Value = 9005.6875
$ws.Range("K31").Value = 2399.5
$ws.Range("L31").Value = 9005.6875
$ws.Range("M31").Value = -2104.5
$ws.Range("N31").Value = -9595.6875
$ws.Range("H34").Value = 7204
$ws.Range("I34").Value = 2399.5
$ws.Range("J34").Value = 9005.6875
$ws.Range("K34").Value = 2399.5
$ws.Range("L34").Value = 9005.6875
$ws.Range("M34").Value = -2197.5
$ws.Range("N34").Value = -9409.6875
$ws.Range("H62").Value = 46772.285
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 46772.285
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 46772.285
$ws.Range("N62").Value = -48020.285
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 46772.285
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 46772.285
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 233861.425
$ws.Range("N65").Value = -240101.425
$ws.Range("M65").ClearContents()
$ws.Range("H99").Value = 8166.4165
$ws.Range("I99").Value = 3999.25
$ws.Range("K99").Value = 3999.25
$ws.Range("M99").Value = -2501.25
$ws.Range("H126").Value = 8166.4165
$ws.Range("I126").Value = 3999.25
$ws.Range("K126").Value = 11997.75
$ws.Range("M126").Value = -9527.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 104.291664
$ws.Range("I2").Value = 94.888885
$ws.Range("J2").Value = 109.933334
$ws.Range("K2").Value = 569.33331
$ws.Range("L2").Value = 659.600004
$ws.Range("M2").Value = -456.33331
$ws.Range("N2").Value = -885.600004
$ws.Range("H12").Value = 208.27272
$ws.Range("J12").Value = 207.125
$ws.Range("L12").Value = 621.375
$ws.Range("N12").Value = -967.375
$ws.Range("H22").Value = 1617.25
$ws.Range("J22").Value = 1409
$ws.Range("L22").Value = 4227
$ws.Range("N22").Value = -4565
$ws.Range("H26").Value = 250191.75
$ws.Range("I26").Value = 333421.66
$ws.Range("J26").Value = 502
$ws.Range("K26").Value = 1000264.98
$ws.Range("L26").Value = 1506
$ws.Range("M26").Value = -999976.98
$ws.Range("N26").Value = -2082
$ws.Range("H27").Value = 1617.25
$ws.Range("J27").Value = 1409
$ws.Range("L27").Value = 4227
$ws.Range("N27").Value = -4431
$ws.Range("H107").Value = 1124.875
$ws.Range("J107").Value = 1161.2
$ws.Range("L107").Value = 3483.6
$ws.Range("N107").Value = -7323.6
$ws.Range("H129").Value = 1912.8572
$ws.Range("J129").Value = 1678
$ws.Range("L129").Value = 5034
$ws.Range("N129").Value = -15034
$ws.Range("H131").Value = 8623380
$ws.Range("I131").Value = 16668950
$ws.Range("J131").Value = 6947219
$ws.Range("K131").Value = 50006850
$ws.Range("L131").Value = 20841657
$ws.Range("M131").Value = -50001810
$ws.Range("N131").Value = -20851737
$ws.Range("H134").Value = 16150
$ws.Range("I134").Value = 5250
$ws.Range("J134").Value = 32500
$ws.Range("K134").Value = 15750
$ws.Range("L134").Value = 97500
$ws.Range("M134").Value = -10680
$ws.Range("N134").Value = -107640

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 128.8077
$ws.Range("I2").Value = 100.73333
$ws.Range("J2").Value = 167.09091
$ws.Range("K2").Value = 100.73333
$ws.Range("L2").Value = 167.09091
$ws.Range("M2").Value = 12.26667
$ws.Range("N2").Value = -393.09091
$ws.Range("H97").Value = 369.5
$ws.Range("I97").Value = 362.25
$ws.Range("K97").Value = 362.25
$ws.Range("M97").Value = 133.75
$ws.Range("H126").Value = 3865.2896
$ws.Range("I126").Value = 3016.4348
$ws.Range("K126").Value = 9049.304400000001
$ws.Range("M126").Value = -6579.304400000001
$ws.Range("H132").Value = 8410.556
$ws.Range("I132").Value = 7956.5713
$ws.Range("K132").Value = 23869.7139
$ws.Range("M132").Value = -21339.7139

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7603
$ws.Range("I46").Value = 7499
$ws.Range("K46").Value = 7499
$ws.Range("M46").Value = -7311
$ws.Range("H68").Value = 5688315.5
$ws.Range("I68").Value = 22727272
$ws.Range("K68").Value = 22727272
$ws.Range("M68").Value = -22726523
$ws.Range("H71").Value = 5688315.5
$ws.Range("I71").Value = 22727272
$ws.Range("K71").Value = 113636360
$ws.Range("M71").Value = -113632616
$ws.Range("H82").Value = 5210198.5
$ws.Range("I82").Value = 7814597.5
$ws.Range("J82").Value = 1400
$ws.Range("K82").Value = 7814597.5
$ws.Range("L82").Value = 1400
$ws.Range("M82").Value = -7814236.5
$ws.Range("N82").Value = -2122
$ws.Range("H85").Value = 5210198.5
$ws.Range("I85").Value = 7814597.5
$ws.Range("J85").Value = 1400
$ws.Range("K85").Value = 7814597.5
$ws.Range("L85").Value = 1400
$ws.Range("M85").Value = -7813349.5
$ws.Range("N85").Value = -3896
$ws.Range("H93").Value = 2581.3333
$ws.Range("I93").Value = 573.5714
$ws.Range("J93").Value = 5392.2
$ws.Range("K93").Value = 573.5714
$ws.Range("L93").Value = 5392.2
$ws.Range("M93").Value = 674.4286
$ws.Range("N93").Value = -7888.2
$ws.Range("H100").Value = 2421.125
$ws.Range("I100").Value = 2421.125
$ws.Range("K100").Value = 2421.125
$ws.Range("M100").Value = -1880.125
$ws.Range("H132").Value = 6575.0713
$ws.Range("I132").Value = 4715.636
$ws.Range("K132").Value = 14146.908
$ws.Range("M132").Value = -11616.908

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 834741.8
$ws.Range("I100").Value = 1177290.4
$ws.Range("K100").Value = 2354580.8
$ws.Range("M100").Value = -2354039.8
$ws.Range("H132").Value = 111114060
$ws.Range("I132").Value = 55555556
$ws.Range("J132").Value = 125003700
$ws.Range("K132").Value = 166666668
$ws.Range("L132").Value = 375011100
$ws.Range("M132").Value = -166666138
$ws.Range("N132").Value = -375016160
